# Apply cell value updates for cryptos.xlsx (Mon Nov 25 06:53:30 UTC 2024 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.132.62"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "3.381.39"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "657.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.45"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.422"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.72%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.04"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("D11").Value = "3.380.37"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.208"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.75"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "97.770.05"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.10"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000255"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.74%  "
$ws.Range("D17").Value = "4.019.51"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "3.381.67"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.98"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.516"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -9.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.99"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "509.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.41"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000200"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.49"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("D29").Value = "3.568.34"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.143"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.59"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.558"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.73"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.90"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "529.18"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.845"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.73"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0425"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.24"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.82"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.27%  "
